$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells stay as text, since many values look numeric
# (e.g. "20.11", "1.00") and would otherwise be auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '51.492.52'
$ws.Range('E2').Value = '  +4.66%  '
$ws.Range('D3').Value = '2.743.13'
$ws.Range('E3').Value = '  +4.22%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '115.89'
$ws.Range('E5').Value = '  +3.76%  '
$ws.Range('D6').Value = '332.46'
$ws.Range('E6').Value = '  +3.00%  '
$ws.Range('E7').Value = '  +2.26%  '
$ws.Range('D8').Value = '1.00'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '0.572'
$ws.Range('E9').Value = '  +5.15%  '
$ws.Range('D10').Value = '41.58'
$ws.Range('E10').Value = '  +4.40%  '
$ws.Range('D11').Value = '0.0855'
$ws.Range('E11').Value = '  +5.41%  '
$ws.Range('D12').Value = '20.11'
$ws.Range('E12').Value = '  +1.47%  '
$ws.Range('E13').Value = '  +2.74%  '
$ws.Range('D14').Value = '7.61'
$ws.Range('E14').Value = '  +4.79%  '
$ws.Range('D15').Value = '3.172.39'
$ws.Range('E15').Value = '  +4.33%  '
$ws.Range('D16').Value = '2.736.57'
$ws.Range('E16').Value = '  +3.44%  '
$ws.Range('D17').Value = '0.881'
$ws.Range('E17').Value = '  +2.40%  '
$ws.Range('D18').Value = '51.468.91'
$ws.Range('E18').Value = '  +4.72%  '
$ws.Range('E19').Value = '  +7.77%  '
$ws.Range('D20').Value = '13.42'
$ws.Range('E20').Value = '  +3.90%  '
$ws.Range('D21').Value = '6.84'
$ws.Range('E21').Value = '  +2.02%  '
$ws.Range('D23').Value = '278.44'
$ws.Range('E23').Value = '  +3.33%  '
$ws.Range('D24').Value = '69.36'
$ws.Range('E24').Value = '  +1.25%  '
$ws.Range('D25').Value = '2.65'
$ws.Range('E25').Value = '  +4.14%  '
$ws.Range('D26').Value = '26.76'
$ws.Range('E26').Value = '  +2.36%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').Value = '10.17'
$ws.Range('E28').Value = '  +1.31%  '
$ws.Range('E29').Value = '  -0.38%  '
$ws.Range('E30').Value = '  +1.97%  '
$ws.Range('D31').Value = '35.03'
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('D32').Value = '49.97'
$ws.Range('E32').Value = '  +0.88%  '
$ws.Range('E33').Value = '  +0.92%  '
$ws.Range('D34').Value = '0.0823'
$ws.Range('E34').Value = '  +3.15%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').Value = '19.12'
$ws.Range('E35').Value = '  +0.44%  '
$ws.Range('B36').Value = 'FirstDigitalUSD'
$ws.Range('C36').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('E38').Value = '  +1.88%  '
$ws.Range('D39').Value = '3.21'
$ws.Range('E39').Value = '  +2.64%  '
$ws.Range('D40').Value = '127.80'
$ws.Range('E40').Value = '  +0.69%  '
$ws.Range('E41').Value = '  +4.38%  '
$ws.Range('D42').Value = '2.29'
$ws.Range('E42').Value = '  +7.80%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0343'
$ws.Range('E43').Value = '  +8.27%  '
$ws.Range('B44').Value = 'Stellar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D44').Value = '0.113'
$ws.Range('E44').Value = '  +2.58%  '
$ws.Range('D45').Value = '2.42'
$ws.Range('E45').Value = '  +13.39%  '
$ws.Range('D46').Value = '2.086.81'
$ws.Range('E46').Value = '  +0.88%  '
$ws.Range('D47').Value = '3.32'
$ws.Range('E47').Value = '  +2.04%  '
$ws.Range('E48').Value = '  +4.25%  '
$ws.Range('E49').Value = '  +6.18%  '
$ws.Range('D50').Value = '8.94'
$ws.Range('E50').Value = '  +0.21%  '
$ws.Range('D51').Value = '59.84'
$ws.Range('E51').Value = '  +2.14%  '
